$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price/volume figures are stored as literal text (e.g. thousand-separated
# "23.118.27" or padded "  -1.12%  "). A leading apostrophe forces Excel to
# keep genuinely numeric-looking strings (like "0.9965") as text instead of
# auto-converting them to a Double, matching the workbook's original
# inline-string cells.
$ws.Range("D2").Value = "23.118.27"
$ws.Range("E2").Value = "  +0.00%  "
$ws.Range("D3").Value = "1.588.66"
$ws.Range("E3").Value = "  -1.12%  "
$ws.Range("D4").Value = "'0.9965"
$ws.Range("E4").Value = "  -0.48%  "
$ws.Range("D5").Value = "'0.9976"
$ws.Range("E5").Value = "  -0.34%  "
$ws.Range("D6").Value = "'301.14"
$ws.Range("E6").Value = "  -0.17%  "
$ws.Range("D7").Value = "'0.3760"
$ws.Range("E7").Value = "  -0.24%  "
$ws.Range("B8").Value = "Cardano"
$ws.Range("C8").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D8").Value = "'0.3603"
$ws.Range("E8").Value = "  -1.18%  "
$ws.Range("B9").Value = "OKB"
$ws.Range("C9").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D9").Value = "'50.68"
$ws.Range("E9").Value = "  +3.91%  "
$ws.Range("D10").Value = "'1.237"
$ws.Range("E10").Value = "  -1.62%  "
$ws.Range("D11").Value = "'0.9961"
$ws.Range("E11").Value = "  -0.50%  "
$ws.Range("D12").Value = "'0.08027"
$ws.Range("E12").Value = "  -0.44%  "
$ws.Range("E13").Value = "  -2.79%  "
$ws.Range("D14").Value = "'6.495"
$ws.Range("E14").Value = "  -1.17%  "
$ws.Range("D15").Value = "'7.355"
$ws.Range("E15").Value = "  -0.29%  "
$ws.Range("D16").Value = "'0.00001240"
$ws.Range("E16").Value = "  -0.83%  "
$ws.Range("D17").Value = "1.588.04"
$ws.Range("E17").Value = "  -1.42%  "
$ws.Range("D18").Value = "'92.86"
$ws.Range("E18").Value = "  +1.74%  "
$ws.Range("D19").Value = "'0.06758"
$ws.Range("E19").Value = "  -0.54%  "
$ws.Range("E20").Value = "  -2.18%  "
$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").Value = "'6.422"
$ws.Range("E21").Value = "  -2.06%  "
$ws.Range("B22").Value = "Dai"
$ws.Range("C22").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D22").Value = "'0.9982"
$ws.Range("E22").Value = "  -0.35%  "
$ws.Range("B23").Value = "Cosmos"
$ws.Range("C23").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D23").Value = "'12.76"
$ws.Range("E23").Value = "  -2.04%  "
$ws.Range("B24").Value = "WrappedBTC"
$ws.Range("C24").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D24").Value = "23.104.90"
$ws.Range("E24").Value = "  -0.11%  "
$ws.Range("D25").Value = "'2.376"
$ws.Range("E25").Value = "  +1.19%  "
$ws.Range("D26").Value = "'2.893"
$ws.Range("E26").Value = "  +4.51%  "
$ws.Range("D27").Value = "'20.83"
$ws.Range("E27").Value = "  -1.09%  "
$ws.Range("D28").Value = "'148.40"
$ws.Range("E28").Value = "  -1.29%  "
$ws.Range("D29").Value = "'5.195"
$ws.Range("E29").Value = "  -0.83%  "
$ws.Range("D30").Value = "'132.31"
$ws.Range("E30").Value = "  -0.29%  "
$ws.Range("D31").Value = "'2.355"
$ws.Range("E31").Value = "  -1.83%  "
$ws.Range("D32").Value = "'6.619"
$ws.Range("E32").Value = "  -2.38%  "
$ws.Range("D33").Value = "1.763.81"
$ws.Range("E33").Value = "  -1.19%  "
$ws.Range("D34").Value = "'0.9525"
$ws.Range("E34").Value = "  -1.05%  "
$ws.Range("D35").Value = "'0.07439"
$ws.Range("E35").Value = "  -2.94%  "
$ws.Range("B36").Value = "FraxShare"
$ws.Range("C36").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D36").Value = "'10.05"
$ws.Range("E36").Value = "  -0.22%  "
$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").Value = "'0.02667"
$ws.Range("E37").Value = "  -3.37%  "
$ws.Range("B38").Value = "Algorand"
$ws.Range("C38").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D38").Value = "'0.2497"
$ws.Range("E38").Value = "  -1.82%  "
$ws.Range("B39").Value = "Stellar"
$ws.Range("C39").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D39").Value = "'0.08777"
$ws.Range("E39").Value = "  -1.51%  "
$ws.Range("D40").Value = "'6.088"
$ws.Range("E40").Value = "  -1.89%  "
$ws.Range("D41").Value = "'1.352"
$ws.Range("E41").Value = "  -2.35%  "
$ws.Range("D42").Value = "'0.7018"
$ws.Range("E42").Value = "  -2.13%  "
$ws.Range("D43").Value = "'12.12"
$ws.Range("E43").Value = "  -5.13%  "
$ws.Range("D44").Value = "'14.88"
$ws.Range("E44").Value = "  -4.99%  "
$ws.Range("D45").Value = "'0.6452"
$ws.Range("E45").Value = "  -2.60%  "
$ws.Range("D46").Value = "'0.9965"
$ws.Range("E46").Value = "  -0.40%  "
$ws.Range("D47").Value = "'3.988"
$ws.Range("E47").Value = "  +0.25%  "
$ws.Range("D48").Value = "'2.268"
$ws.Range("E48").Value = "  -0.98%  "
$ws.Range("D49").Value = "'131.15"
$ws.Range("E49").Value = "  +0.59%  "
$ws.Range("D50").Value = "'0.07882"
$ws.Range("E50").Value = "  -1.18%  "
$ws.Range("D51").Value = "'1.203"
$ws.Range("E51").Value = "  +2.71%  "
